$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 89
Set-TextCell "A89" "6399"
Set-TextCell "B89" "7/14/2025"
Set-TextCell "C89" "ESCALADA AV. 966"
Set-TextCell "D89" "9"
Set-TextCell "E89" "808258198"
Set-TextCell "F89" "AYKO"
Set-TextCell "G89" "Pendiente"
Set-TextCell "H89" "Picada"
$ws.Range("I89").Value = 1
Set-TextCell "J89" "Cambio"
Set-TextCell "K89" "Sin equipos"
Set-TextCell "L89" "Pasante"
$ws.Range("M89").Value = -58.493069
$ws.Range("N89").Value = -34.646557
Set-TextCell "O89" "Devoto"
Set-TextCell "P89" "Capital Norte"

# Row 90
Set-TextCell "A90" "-515"
Set-TextCell "B90" "7/15/2025"
Set-TextCell "C90" "Rivadavia 7470"
Set-TextCell "D90" "7"
Set-TextCell "E90" "808263485"
Set-TextCell "F90" "AYKO"
Set-TextCell "G90" "Pendiente"
Set-TextCell "H90" "Picada"
$ws.Range("I90").Value = 1
Set-TextCell "J90" "Cambio"
Set-TextCell "K90" "Sin equipos"
Set-TextCell "L90" "Pasante"
$ws.Range("M90").Value = -58.470715
$ws.Range("N90").Value = -34.631107
Set-TextCell "O90" "Boedo"
Set-TextCell "P90" "Capital Sur"

# Row 91
Set-TextCell "A91" "6398"
Set-TextCell "B91" "7/15/2025"
Set-TextCell "C91" "LARRAZABAL AV. 579"
Set-TextCell "D91" "9"
Set-TextCell "E91" "808373655"
Set-TextCell "F91" "AYKO"
Set-TextCell "G91" "Pendiente"
Set-TextCell "H91" "Poste inclinado"
$ws.Range("I91").Value = 1
Set-TextCell "J91" "Aplomo"
Set-TextCell "K91" "Sin equipos"
Set-TextCell "L91" "Poste"
$ws.Range("M91").Value = -58.510247
$ws.Range("N91").Value = -34.645038
Set-TextCell "O91" "Devoto"
Set-TextCell "P91" "Capital Norte"

# Row 92
Set-TextCell "A92" "6410"
Set-TextCell "B92" "7/15/2025"
Set-TextCell "C92" "BEAUCHEF 421"
Set-TextCell "D92" "6"
Set-TextCell "E92" "808373658"
Set-TextCell "F92" "AYKO"
Set-TextCell "G92" "Pendiente"
Set-TextCell "H92" "Picada"
$ws.Range("I92").Value = 1
Set-TextCell "J92" "Cambio"
Set-TextCell "K92" "Sin equipos"
Set-TextCell "L92" "Pasante"
$ws.Range("M92").Value = -58.433724
$ws.Range("N92").Value = -34.621643
Set-TextCell "O92" "Boedo"
Set-TextCell "P92" "Capital Sur"

# Row 93
Set-TextCell "A93" "-518"
Set-TextCell "B93" "7/16/2025"
Set-TextCell "C93" "Crisologo Larralde 4073"
Set-TextCell "D93" "12"
Set-TextCell "E93" "808373627"
Set-TextCell "F93" "AYKO"
Set-TextCell "G93" "Pendiente"
Set-TextCell "H93" "Cambiar columna corroída en base"
$ws.Range("I93").Value = 1
Set-TextCell "J93" "Cambio"
Set-TextCell "K93" "Sin equipos"
Set-TextCell "L93" "Pasante"
$ws.Range("M93").Value = -58.483145
$ws.Range("N93").Value = -34.557043
Set-TextCell "O93" "Saavedra"
Set-TextCell "P93" "Capital Norte"
